# MCT-1A schedule fix: "6 hours by turn fix"
# Shifts several time-slot rows and rebalances which weekday/period a
# couple of classes (João Rodrigues / José Ferreira) fall on, and inserts
# an extra half-hour slot at the end of the day (adds rows 15-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (08:40 slot): move "José Ferreira" class from D4 to E4, clear C4/D4
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "José Ferreira-Tecnologia dos Materiais"

# --- Row 6 (09:50 slot): clear E6 (José Ferreira now only in D6)
$ws.Range("E6").Value = "-"

# --- Row 7 (10:40 slot): add João Rodrigues / José Ferreira classes
$ws.Range("C7").Value = "João Rodrigues-Desenho Técnico"
$ws.Range("D7").Value = "José Ferreira-Tecnologia dos Materiais"

# --- Row 8 (11:30 slot): no longer lunch
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"

# --- Row 9: becomes the lunch ("Almoço") slot at 12:20
$ws.Range("A9").Value = "12:20"
$ws.Range("B9").Value = "Almoço"
$ws.Range("C9").Value = "Almoço"
$ws.Range("D9").Value = "Almoço"
$ws.Range("E9").Value = "Almoço"
$ws.Range("F9").Value = "Almoço"

# --- Row 10: time shifts earlier (content unchanged)
$ws.Range("A10").Value = "13:00"

# --- Row 11: time shifts earlier (content unchanged)
$ws.Range("A11").Value = "13:50"

# --- Row 12: time shifts earlier, no longer "Intervalo"
$ws.Range("A12").Value = "14:40"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"

# --- Row 13: time shifts earlier, becomes "Intervalo"
$ws.Range("A13").Value = "15:30"
$ws.Range("B13").Value = "Intervalo"
$ws.Range("C13").Value = "Intervalo"
$ws.Range("D13").Value = "Intervalo"
$ws.Range("E13").Value = "Intervalo"
$ws.Range("F13").Value = "Intervalo"

# --- Row 14: time shifts earlier (content stays "-")
$ws.Range("A14").Value = "15:50"
$ws.Range("B14").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"

# --- New row 15: what used to be the old row 14 (16:40 slot)
$ws.Range("A15").Value = "16:40"
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"

# --- New row 16: brand-new 17:30 slot
$ws.Range("A16").Value = "17:30"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"

# --- New row 17: brand-new 18:20 slot, label only (blank turns).
# A leading apostrophe forces Excel to store an explicit empty *text*
# cell (instead of clearing/removing it entirely), matching the template
# row's blank-but-present cells; reset the style afterwards so the
# quote-prefix formatting doesn't stick.
$ws.Range("A17").Value = "18:20"
$ws.Range("B17").Value = "'"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = "'"
$ws.Range("F17").Style = "Normal"
